$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Containers")

$ws.Range("A9").Value = "Admin"
$ws.Range("B9").Value = "messagestodeveloper"
$ws.Range("C9").Value = "Contains a series of messages to owner."
$ws.Range("D9").Value = "YYYY.MM.DD.HH.MM.SS"

$ws.Range("D9").Select()
